$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1246.4117
$ws.Range("J19").Value = 1300.2142
$ws.Range("L19").Value = 1300.2142
$ws.Range("N19").Value = -1650.2142
$ws.Range("H28").Value = 1558.875
$ws.Range("I28").Value = 1460.8889
$ws.Range("J28").Value = 1684.8572
$ws.Range("K28").Value = 1460.8889
$ws.Range("L28").Value = 1684.8572
$ws.Range("M28").Value = -975.8888999999999
$ws.Range("N28").Value = -2654.8572
$ws.Range("H31").Value = 200
$ws.Range("I31").Value = 200
$ws.Range("K31").Value = 600
$ws.Range("M31").Value = -370
$ws.Range("H138").Value = 3875.5789
$ws.Range("J138").Value = 4664.442
$ws.Range("L138").Value = 13993.326
$ws.Range("N138").Value = -24273.326

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10240.786
$ws.Range("I45").Value = 13615.889
$ws.Range("J45").Value = 4165.6
$ws.Range("K45").Value = 13615.889
$ws.Range("L45").Value = 4165.6
$ws.Range("M45").Value = -13238.889
$ws.Range("N45").Value = -4919.6
$ws.Range("H110").Value = 11520.042
$ws.Range("I110").Value = 17907.363
$ws.Range("J110").Value = 6115.385
$ws.Range("K110").Value = 17907.363
$ws.Range("L110").Value = 6115.385
$ws.Range("M110").Value = -15862.363
$ws.Range("N110").Value = -10205.385
$ws.Range("H134").Value = 53250
$ws.Range("J134").Value = 53250
$ws.Range("L134").Value = 53250
$ws.Range("N134").Value = -63390

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 197414.48
$ws.Range("I22").Value = 310.14285
$ws.Range("J22").Value = 342649.25
$ws.Range("K22").Value = 310.14285
$ws.Range("L22").Value = 342649.25
$ws.Range("M22").Value = -137.14285
$ws.Range("N22").Value = -342995.25
$ws.Range("H134").Value = 2264.122
$ws.Range("I134").Value = 2014.2972
$ws.Range("K134").Value = 6042.8916
$ws.Range("M134").Value = -3507.8916
$ws.Range("H138").Value = 99989
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 469.55554
$ws.Range("I22").Value = 245
$ws.Range("J22").Value = 533.7143
$ws.Range("K22").Value = 245
$ws.Range("L22").Value = 533.7143
$ws.Range("M22").Value = 105
$ws.Range("N22").Value = -1233.7143
$ws.Range("H26").Value = 14509.5
$ws.Range("I26").Value = 9019
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 9019
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = -8732
$ws.Range("N26").Value = -20574
$ws.Range("H31").Value = 33874.605
$ws.Range("I31").Value = 51607.65
$ws.Range("K31").Value = 51607.65
$ws.Range("M31").Value = -51312.65
$ws.Range("H34").Value = 33874.605
$ws.Range("I34").Value = 51607.65
$ws.Range("K34").Value = 51607.65
$ws.Range("M34").Value = -51405.65
$ws.Range("H132").Value = 4941.75
$ws.Range("I132").Value = 4896.4736
$ws.Range("J132").Value = 5113.8
$ws.Range("K132").Value = 14689.4208
$ws.Range("L132").Value = 15341.4
$ws.Range("M132").Value = -12159.4208
$ws.Range("N132").Value = -20401.4
$ws.Range("H134").Value = 12376.441
$ws.Range("I134").Value = 6748.3887
$ws.Range("K134").Value = 20245.1661
$ws.Range("M134").Value = -17710.1661
$ws.Range("H137").Value = 57942.668
$ws.Range("J137").Value = 54560.5
$ws.Range("L137").Value = 54560.5
$ws.Range("N137").Value = -64760.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 75.71429000000001
$ws.Range("I2").Value = 96.111115
$ws.Range("J2").Value = 39
$ws.Range("K2").Value = 576.66669
$ws.Range("L2").Value = 234
$ws.Range("M2").Value = -463.66669
$ws.Range("N2").Value = -460
$ws.Range("H58").Value = 1000
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H113").Value = 396.85715
$ws.Range("J113").Value = 364.625
$ws.Range("L113").Value = 1093.875
$ws.Range("N113").Value = -5433.875
$ws.Range("H121").Value = 7542
$ws.Range("I121").Value = 4109.9
$ws.Range("K121").Value = 12329.7
$ws.Range("M121").Value = -11019.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 260.81818
$ws.Range("I2").Value = 272.38095
$ws.Range("K2").Value = 272.38095
$ws.Range("M2").Value = -159.38095
$ws.Range("H102").Value = 1642.2273
$ws.Range("I102").Value = 1654
$ws.Range("K102").Value = 1654
$ws.Range("M102").Value = -32
$ws.Range("H132").Value = 1432785.8
$ws.Range("I132").Value = 2003500
$ws.Range("K132").Value = 6010500
$ws.Range("M132").Value = -6007970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1575
$ws.Range("I22").Value = 725
$ws.Range("J22").Value = 1745
$ws.Range("K22").Value = 725
$ws.Range("L22").Value = 1745
$ws.Range("M22").Value = -430
$ws.Range("N22").Value = -2335
$ws.Range("H27").Value = 1575
$ws.Range("I27").Value = 725
$ws.Range("J27").Value = 1745
$ws.Range("K27").Value = 725
$ws.Range("L27").Value = 1745
$ws.Range("M27").Value = -618
$ws.Range("N27").Value = -1959

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 16327.143
$ws.Range("I28").Value = 10000
$ws.Range("K28").Value = 10000
$ws.Range("M28").Value = -9652
$ws.Range("H74").Value = 12066.667
$ws.Range("J74").Value = 12066.667
$ws.Range("L74").Value = 12066.667
$ws.Range("N74").Value = -13938.667
$ws.Range("H77").Value = 12066.667
$ws.Range("J77").Value = 12066.667
$ws.Range("L77").Value = 36200.001
$ws.Range("N77").Value = -45560.001
$ws.Range("H126").Value = 3931
$ws.Range("I126").Value = 3755.5
$ws.Range("J126").Value = 4750
$ws.Range("K126").Value = 11266.5
$ws.Range("L126").Value = 14250
$ws.Range("M126").Value = -8796.5
$ws.Range("N126").Value = -19190
$ws.Range("H132").Value = 2815.2334
$ws.Range("J132").Value = 2324.111
$ws.Range("L132").Value = 6972.333
$ws.Range("N132").Value = -12032.333
$ws.Range("H137").Value = 125894
$ws.Range("J137").Value = 125894
$ws.Range("L137").Value = 125894
$ws.Range("N137").Value = -136094
$ws.Range("H141").Value = 165125
$ws.Range("J141").Value = 165125
$ws.Range("L141").Value = 165125
$ws.Range("N141").Value = -175485
